# Optimization of projections variables
# Update the three player-stack tables on Sheet1 with new team/player names
# and new salary/FPTS figures, per the latest projections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Table 1 (columns A:D) -> New York Mets hitters ----
$ws.Range("A1").Value2 = "New York Mets hitters"

$ws.Range("A3").Value2 = "Conforto"
$ws.Range("B3").Value2 = 3100
$ws.Range("C3").Value2 = 21.7

$ws.Range("A4").Value2 = "Cespedes"
$ws.Range("B4").Value2 = 4000
$ws.Range("C4").Value2 = 18.7

$ws.Range("A5").Value2 = "Bruce"
$ws.Range("B5").Value2 = 3100
$ws.Range("C5").Value2 = 31.2

$ws.Range("A6").Value2 = "Gonzalez"
$ws.Range("B6").Value2 = 2500
$ws.Range("C6").Value2 = 43.6

$ws.Range("D8").Value2 = "Success"

# ---- Table 2 (columns F:I) -> Detroit Tigers hitters ----
$ws.Range("F1").Value2 = "Detroit Tigers hitters"

$ws.Range("F3").Value2 = "Martin"
$ws.Range("H3").Value2 = 34.4

$ws.Range("F4").Value2 = "Martinez"
$ws.Range("G4").Value2 = 2600
$ws.Range("H4").Value2 = 10

$ws.Range("F5").Value2 = "Hicks"
$ws.Range("G5").Value2 = 2600
$ws.Range("H5").Value2 = 3

$ws.Range("F6").Value2 = "McCann"
$ws.Range("G6").Value2 = 2400
$ws.Range("H6").Value2 = 3

# ---- Table 3 (columns K:N) -> San Francisco Giants lefties ----
$ws.Range("K1").Value2 = "San Francisco Giants lefties"

$ws.Range("K3").Value2 = "Blanco"
$ws.Range("L3").Value2 = 2400
$ws.Range("M3").Value2 = 3

$ws.Range("K4").Value2 = "Belt"
$ws.Range("L4").Value2 = 3800
$ws.Range("M4").Value2 = 3

$ws.Range("K5").Value2 = "Crawford"
$ws.Range("L5").Value2 = 2800
$ws.Range("M5").Value2 = 3

$ws.Range("K6").Value2 = "Hanson"
$ws.Range("L6").Value2 = 3100
$ws.Range("M6").Value2 = 6

$ws.Range("N8").Value2 = "Failure"

# Recalculate everything (also drives Sheet2 formula results)
$excel.CalculateFull()

# Restore the active cell selection recorded for Sheet1
$ws.Range("K13").Select()
